$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(303, 1562.231466726676),
    @(304, 3149.240758891025),
    @(305, 4761.455215897783),
    @(306, 6399.307790093049),
    @(307, 8063.236924101251),
    @(308, 9753.68641156211),
    @(309, 11471.10524059318),
    @(310, 13215.9474193051),
    @(311, 14988.67178270165),
    @(312, 16789.74178030778),
    @(313, 18619.6252438775),
    @(314, 20478.79413454716),
    @(315, 22367.72426881849),
    @(316, 24286.89502277552),
    @(317, 26236.78901396181),
    @(318, 28217.89176037541),
    @(319, 30230.69131606622),
    @(320, 32275.67788285805),
    @(321, 34353.34339775828),
    @(322, 36464.18109565748),
    @(323, 38608.68504697403),
    @(324, 40787.34966994636),
    @(325, 43000.66921733767),
    @(326, 45249.13723737408),
    @(327, 47533.24600880643),
    @(328, 49853.48595005651),
    @(329, 52210.34500248212),
    @(330, 54604.30798787673),
    @(331, 57035.85594040353),
    @(332, 59505.46541325291),
    @(333, 62013.60776040423),
    @(334, 64560.74839397207),
    @(335, 67147.3460177168),
    @(336, 69773.8518374049),
    @(337, 72440.70874881349),
    @(338, 75148.35050428357),
    @(339, 77897.20085884085),
    @(340, 80687.6726970223),
    @(341, 83520.16714166098),
    @(342, 86395.07264600335),
    @(343, 89312.76407065538),
    @(344, 92273.60174697166),
    @(345, 95277.93052862627),
    @(346, 98326.07883321807),
    @(347, 101418.3576758894),
    @(348, 104555.0596970436),
    @(349, 107736.4581863681),
    @(350, 110962.806105471),
    @(351, 114234.33511155),
    @(352, 117551.2545846056),
    @(353, 120913.7506608076),
    @(354, 124321.9852747057),
    @(355, 127776.0952130563),
    @(356, 131276.1911831023),
    @(357, 134822.3568982082),
    @(358, 138905.4527780479),
    @(359, 143041.2135274808),
    @(360, 147229.6519363095),
    @(361, 151470.7470686474),
    @(362, 155764.4434059678),
    @(363, 160110.6500308372),
    @(364, 164509.2398546874),
    @(365, 168960.0488929176),
    @(366, 173462.8755905561),
    @(367, 178017.4802016286),
    @(368, 182623.5842252702),
    @(369, 187280.8699015069),
    @(370, 191988.9797694895),
    @(371, 196747.5162908118),
    @(372, 201556.0415403774),
    @(373, 205654.4506582586),
    @(374, 209793.6205756541),
    @(375, 213973.0403188611),
    @(376, 218192.1626352042),
    @(377, 222450.4039863175),
    @(378, 226747.1446064895),
    @(379, 231081.7286268641),
    @(380, 235453.4642660647),
    @(381, 239861.6240875848),
    @(382, 244305.4453240507),
    @(383, 247542.0213823364),
    @(384, 250801.4648262378),
    @(385, 254083.0435117222),
    @(386, 257385.9968319098),
    @(387, 260709.5360696641),
    @(388, 264052.8447942914),
    @(389, 267415.0793014182),
    @(390, 270795.3690949908),
    @(391, 274192.8174102352),
    @(392, 277606.5017763217),
    @(393, 283307.9551854341),
    @(394, 289039.9962735071),
    @(395, 294801.3477246399),
    @(396, 300590.7073699429),
    @(397, 306406.7497169889),
    @(398, 312248.1275298574),
    @(399, 318113.4734563555),
    @(400, 324001.4016989213),
    @(401, 329910.5097256748),
    @(402, 335839.3800180302),
    @(403, 351859.6762067491),
    @(404, 367952.9128960276),
    @(405, 384116.5541667915),
    @(406, 400348.0852220349),
    @(407, 416645.0180621215),
    @(408, 433004.8971705797),
    @(409, 449425.3051955232),
    @(410, 465903.8686113079),
    @(411, 482438.2633446286),
    @(412, 499026.2203486587),
    @(413, 514847.7654637888),
    @(414, 530714.8940952768),
    @(415, 546625.5831641678),
    @(416, 562577.8809572249),
    @(417, 578569.9119169221),
    @(418, 594599.8812759762),
    @(419, 610666.0795155464),
    @(420, 626766.8866251045),
    @(421, 642900.7761407595),
    @(422, 659066.3189374893),
    @(423, 675262.1867493393),
    @(424, 691487.1553901108),
    @(425, 707740.1076455261),
    @(426, 724020.035806153),
    @(427, 740326.0438087285),
    @(428, 756657.3489517822),
    @(429, 773013.2831497326),
    @(430, 789393.2936879841),
    @(431, 805796.9434399192),
    @(432, 822223.9105052016),
    @(433, 838673.9872275017),
    @(434, 855147.0785486251),
    @(435, 871643.1996552453),
    @(436, 888162.4728739177),
    @(437, 904705.1237700412),
    @(438, 921271.476406802),
    @(439, 937861.9477211658),
    @(440, 954477.0409755436),
    @(441, 971117.3382461244),
    @(442, 987783.4919119548),
    @(443, 1004476.215112791),
    @(444, 1021196.27114867),
    @(445, 1037944.461799973),
    @(446, 1054721.614553673),
    @(447, 1071528.568729396),
    @(448, 1088366.160507993),
    @(449, 1105235.206875488),
    @(450, 1122136.488506451),
    @(451, 1139070.731623216),
    @(452, 1156038.588880617)
)

foreach ($pair in $values) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 2).Value = $v
}

Write-Output "Updated $($values.Count) cells in column B (rows 303-452)"